$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 with the "Fragen" entry about paginating to 20 questions per page
$ws.Range("A17").Value = "Fragen"
$ws.Range("B17").Value = "Zum Blättern in Seiten unterteilen, 20 Fragen pro Seite"
$ws.Range("C17").Value = "HS"

# Copy formatting (centered alignment) from the row above, matching other "Wer?" cells
[void]$ws.Range("C16").Copy()
[void]$ws.Range("C17").PasteSpecial(-4122)

# Move the active selection down to A18, as it was after inserting this row
[void]$ws.Range("A18").Select()
